$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values
$ws.Range("B2").Value = 11.434701127519407
$ws.Range("C2").Value = 11.142360010331798
$ws.Range("D2").Value = 12.478975489976495
$ws.Range("E2").Value = 11.303159915002924

# Row 3 data values
$ws.Range("B3").Value = 10.818102188479116
$ws.Range("C3").Value = 9.2517744904732471
$ws.Range("D3").Value = 10.789385201668139
$ws.Range("E3").Value = 10.937539383034414

# Update the selection to reflect the narrower range used
$ws.Range("B1:E3").Select()
